# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Latest Handoff / HO Xliff Generate timestamps bump forward a couple
#   of minutes to reflect the new handoff generation run
# - Status/date columns widen (autofit) to accommodate the longer text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bump forward ---
$overview.Range("G2").Value = "2016-09-06 15:23:15"
$dede.Range("H2").Value     = "2016-09-06 15:23:15"
$zhcn.Range("H2").Value     = "2016-09-06 15:23:00"

# --- Column widths widen to fit the new, longer status text ---
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth     = 16.3
$dede.Columns.Item(3).ColumnWidth     = 16.3
